$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 500.89474
$ws.Range("I2").Value = 191.58333
$ws.Range("K2").Value = 191.58333
$ws.Range("M2").Value = -78.58332999999999
$ws.Range("H19").Value = 1186.8667
$ws.Range("J19").Value = 1023
$ws.Range("L19").Value = 1023
$ws.Range("N19").Value = -1373
$ws.Range("H51").Value = 19799.8
$ws.Range("J51").Value = 19799.8
$ws.Range("L51").Value = 19799.8
$ws.Range("N51").Value = -20767.8
$ws.Range("H69").Value = 12534.056
$ws.Range("J69").Value = 12320.929
$ws.Range("L69").Value = 36962.787
$ws.Range("N69").Value = -38710.787
$ws.Range("H72").Value = 12534.056
$ws.Range("J72").Value = 12320.929
$ws.Range("L72").Value = 110888.361
$ws.Range("N72").Value = -119624.361
$ws.Range("H94").Value = 1113.7778
$ws.Range("I94").Value = 1113.7778
$ws.Range("K94").Value = 1113.7778
$ws.Range("M94").Value = -662.7778000000001
$ws.Range("H137").Value = 1318.3077
$ws.Range("I137").Value = 1281.3
$ws.Range("J137").Value = 1441.6666
$ws.Range("K137").Value = 3843.9
$ws.Range("L137").Value = 4324.9998
$ws.Range("M137").Value = -1293.9
$ws.Range("N137").Value = -9424.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3677371
$ws.Range("I2").Value = 4902361.5
$ws.Range("K2").Value = 4902361.5
$ws.Range("M2").Value = -4902248.5
$ws.Range("H97").Value = 827.4761999999999
$ws.Range("I97").Value = 806.7646999999999
$ws.Range("J97").Value = 915.5
$ws.Range("K97").Value = 806.7646999999999
$ws.Range("L97").Value = 915.5
$ws.Range("M97").Value = -310.7646999999999
$ws.Range("N97").Value = -1907.5
$ws.Range("H102").Value = 8334635
$ws.Range("I102").Value = 10001281
$ws.Range("K102").Value = 10001281
$ws.Range("M102").Value = -9999659
$ws.Range("H116").Value = 3677371
$ws.Range("I116").Value = 4902361.5
$ws.Range("K116").Value = 4902361.5
$ws.Range("M116").Value = -4900067.5
$ws.Range("H132").Value = 1669821.8
$ws.Range("I132").Value = 2043541.2
$ws.Range("K132").Value = 6130623.6
$ws.Range("M132").Value = -6128093.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3677371
$ws.Range("I3").Value = 4902361.5
$ws.Range("K3").Value = 4902361.5
$ws.Range("M3").Value = -4902247.5
$ws.Range("H86").Value = 2699.7666
$ws.Range("J86").Value = 2398.9167
$ws.Range("L86").Value = 2398.9167
$ws.Range("N86").Value = -4644.9167
$ws.Range("H89").Value = 2699.7666
$ws.Range("J89").Value = 2398.9167
$ws.Range("L89").Value = 11994.5835
$ws.Range("N89").Value = -23226.5835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13200.762
$ws.Range("I31").Value = 4246.222
$ws.Range("K31").Value = 4246.222
$ws.Range("M31").Value = -3951.222
$ws.Range("H34").Value = 13200.762
$ws.Range("I34").Value = 4246.222
$ws.Range("K34").Value = 4246.222
$ws.Range("M34").Value = -4044.222
$ws.Range("H97").Value = 39990
$ws.Range("J97").Value = 39990
$ws.Range("L97").Value = 39990
$ws.Range("N97").Value = -41972
$ws.Range("H107").Value = 1189964.1
$ws.Range("I107").Value = 1812210
$ws.Range("J107").Value = 256595.25
$ws.Range("K107").Value = 1812210
$ws.Range("L107").Value = 256595.25
$ws.Range("M107").Value = -1810290
$ws.Range("N107").Value = -260435.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3979.4119
$ws.Range("I131").Value = 3998.3333
$ws.Range("J131").Value = 3969.0908
$ws.Range("K131").Value = 11994.9999
$ws.Range("L131").Value = 11907.2724
$ws.Range("M131").Value = -6954.999899999999
$ws.Range("N131").Value = -21987.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 483.92856
$ws.Range("I2").Value = 550.3043
$ws.Range("K2").Value = 550.3043
$ws.Range("M2").Value = -437.3043
$ws.Range("H18").Value = 48999
$ws.Range("I18").Value = 48999
$ws.Range("J18").Value = 48999
$ws.Range("K18").Value = 48999
$ws.Range("L18").Value = 48999
$ws.Range("M18").Value = -48706
$ws.Range("N18").Value = -49585
$ws.Range("H80").Value = 2749.5
$ws.Range("I80").Value = 2749.5
$ws.Range("K80").Value = 2749.5
$ws.Range("M80").Value = -1751.5
$ws.Range("H83").Value = 2749.5
$ws.Range("I83").Value = 2749.5
$ws.Range("K83").Value = 13747.5
$ws.Range("M83").Value = -8755.5
$ws.Range("H97").Value = 1040.3684
$ws.Range("I97").Value = 633.25
$ws.Range("K97").Value = 633.25
$ws.Range("M97").Value = -137.25
$ws.Range("H132").Value = 6948699
$ws.Range("I132").Value = 7357152
$ws.Range("K132").Value = 22071456
$ws.Range("M132").Value = -22068926

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2074.5293
$ws.Range("I16").Value = 861.1818
$ws.Range("J16").Value = 4299
$ws.Range("K16").Value = 861.1818
$ws.Range("L16").Value = 4299
$ws.Range("M16").Value = -691.1818
$ws.Range("N16").Value = -4639
$ws.Range("H40").Value = 2631.3333
$ws.Range("I40").Value = 2660.25
$ws.Range("K40").Value = 2660.25
$ws.Range("M40").Value = -2524.25
$ws.Range("H46").Value = 3642.2856
$ws.Range("I46").Value = 1499.2
$ws.Range("K46").Value = 1499.2
$ws.Range("M46").Value = -1311.2
$ws.Range("H61").Value = 3741.8572
$ws.Range("I61").Value = 3532.1667
$ws.Range("K61").Value = 3532.1667
$ws.Range("M61").Value = -3330.1667
$ws.Range("H113").Value = 3741.8572
$ws.Range("I113").Value = 3532.1667
$ws.Range("K113").Value = 3532.1667
$ws.Range("M113").Value = -1362.1667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 39389
$ws.Range("J93").Value = 39389
$ws.Range("L93").Value = 39389
$ws.Range("N93").Value = -44381
$ws.Range("H132").Value = 15631029
$ws.Range("I132").Value = 23812536
$ws.Range("J132").Value = 11789.546
$ws.Range("K132").Value = 71437608
$ws.Range("L132").Value = 35368.638
$ws.Range("M132").Value = -71435078
$ws.Range("N132").Value = -40428.638
$ws.Range("H136").Value = 10419121
$ws.Range("I136").Value = 11113593
$ws.Range("J136").Value = 2040.6666
$ws.Range("K136").Value = 33340779
$ws.Range("L136").Value = 6121.9998
$ws.Range("M136").Value = -33338229
$ws.Range("N136").Value = -11221.9998
